$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1055.5714
$ws.Range("I19").Value = 500.5
$ws.Range("J19").Value = 1277.6
$ws.Range("K19").Value = 500.5
$ws.Range("L19").Value = 1277.6
$ws.Range("M19").Value = -325.5
$ws.Range("N19").Value = -1627.6
$ws.Range("H80").Value = 134309.6
$ws.Range("I80").Value = 167597.92
$ws.Range("J80").Value = 1156.3334
$ws.Range("K80").Value = 502793.76
$ws.Range("L80").Value = 3469.0002
$ws.Range("M80").Value = -501795.76
$ws.Range("N80").Value = -5465.0002
$ws.Range("H83").Value = 134309.6
$ws.Range("I83").Value = 167597.92
$ws.Range("J83").Value = 1156.3334
$ws.Range("K83").Value = 1508381.28
$ws.Range("L83").Value = 10407.0006
$ws.Range("M83").Value = -1503389.28
$ws.Range("N83").Value = -20391.0006
$ws.Range("H98").Value = 121100
$ws.Range("I98").Value = 138875
$ws.Range("J98").Value = 50000
$ws.Range("K98").Value = 138875
$ws.Range("L98").Value = 50000
$ws.Range("M98").Value = -137377
$ws.Range("N98").Value = -52996
$ws.Range("H106").Value = 8020.5713
$ws.Range("J106").Value = 6373.75
$ws.Range("L106").Value = 6373.75
$ws.Range("N106").Value = -7635.75
$ws.Range("H111").Value = 1629.8334
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1629.8334
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 4889.5002
$ws.Range("N111").Value = -11023.5002
$ws.Range("M111").ClearContents()
$ws.Range("H113").Value = 18675.584
$ws.Range("J113").Value = 8375
$ws.Range("L113").Value = 8375
$ws.Range("N113").Value = -14883
$ws.Range("H122").Value = 121100
$ws.Range("I122").Value = 138875
$ws.Range("J122").Value = 50000
$ws.Range("K122").Value = 416625
$ws.Range("L122").Value = 150000
$ws.Range("M122").Value = -414175
$ws.Range("N122").Value = -154900
$ws.Range("H137").Value = 4233161.5
$ws.Range("I137").Value = 688292.3
$ws.Range("J137").Value = 7360987
$ws.Range("K137").Value = 2064876.9
$ws.Range("L137").Value = 22082961
$ws.Range("M137").Value = -2062326.9
$ws.Range("N137").Value = -22088061
$ws.Range("H138").Value = 8820.121
$ws.Range("I138").Value = 8388.833
$ws.Range("J138").Value = 9066.571
$ws.Range("K138").Value = 25166.499
$ws.Range("L138").Value = 27199.713
$ws.Range("M138").Value = -20026.499
$ws.Range("N138").Value = -37479.713
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 6771.5
$ws.Range("I141").Value = 5157.75
$ws.Range("J141").Value = 9999
$ws.Range("K141").Value = 15473.25
$ws.Range("L141").Value = 29997
$ws.Range("M141").Value = -10293.25
$ws.Range("N141").Value = -40357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2465.3076
$ws.Range("I32").Value = 2371.551
$ws.Range("K32").Value = 2371.551
$ws.Range("M32").Value = -2084.551
$ws.Range("H44").Value = 84988.5
$ws.Range("J44").Value = 84988.5
$ws.Range("L44").Value = 84988.5
$ws.Range("N44").Value = -85964.5
$ws.Range("H61").Value = 12895.23
$ws.Range("I61").Value = 18805.572
$ws.Range("K61").Value = 18805.572
$ws.Range("M61").Value = -18593.572
$ws.Range("H97").Value = 6604.1113
$ws.Range("I97").Value = 7444.933
$ws.Range("K97").Value = 7444.933
$ws.Range("M97").Value = -6948.933
$ws.Range("H132").Value = 4646.2
$ws.Range("I132").Value = 4374.5835
$ws.Range("K132").Value = 13123.7505
$ws.Range("M132").Value = -10593.7505
$ws.Range("H136").Value = 12895.23
$ws.Range("I136").Value = 18805.572
$ws.Range("K136").Value = 56416.716
$ws.Range("M136").Value = -53866.716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 72996.2
$ws.Range("J35").Value = 72996.2
$ws.Range("L35").Value = 72996.2
$ws.Range("N35").Value = -73616.2
$ws.Range("H94").Value = 8280.146
$ws.Range("I94").Value = 10440.311
$ws.Range("K94").Value = 10440.311
$ws.Range("M94").Value = -9989.311
$ws.Range("H105").Value = 83428.5
$ws.Range("I105").Value = 141126.75
$ws.Range("K105").Value = 141126.75
$ws.Range("M105").Value = -139379.75
$ws.Range("H134").Value = 10457.158
$ws.Range("I134").Value = 11402.25
$ws.Range("K134").Value = 34206.75
$ws.Range("M134").Value = -31671.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2876.6511
$ws.Range("I31").Value = 2160.8667
$ws.Range("J31").Value = 3260.1072
$ws.Range("K31").Value = 2160.8667
$ws.Range("L31").Value = 3260.1072
$ws.Range("M31").Value = -1865.8667
$ws.Range("N31").Value = -3850.1072
$ws.Range("H34").Value = 2876.6511
$ws.Range("I34").Value = 2160.8667
$ws.Range("J34").Value = 3260.1072
$ws.Range("K34").Value = 2160.8667
$ws.Range("L34").Value = 3260.1072
$ws.Range("M34").Value = -1958.8667
$ws.Range("N34").Value = -3664.1072
$ws.Range("H132").Value = 20866928
$ws.Range("I132").Value = 25651382
$ws.Range("K132").Value = 76954146
$ws.Range("M132").Value = -76951616

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 240.8125
$ws.Range("I2").Value = 184.85715
$ws.Range("K2").Value = 1109.1429
$ws.Range("M2").Value = -996.1428999999998
$ws.Range("H12").Value = 9.363636
$ws.Range("I12").Value = 17
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 51
$ws.Range("L12").Value = 15
$ws.Range("M12").Value = 122
$ws.Range("N12").Value = -361
$ws.Range("H88").Value = 4714.2856
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15856
$ws.Range("H91").Value = 4714.2856
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 15000
$ws.Range("N91").Value = -17964
$ws.Range("H92").Value = 560.8333
$ws.Range("I92").Value = 483.85715
$ws.Range("K92").Value = 1451.57145
$ws.Range("M92").Value = -203.5714499999999
$ws.Range("H99").Value = 6791
$ws.Range("I99").Value = 1992.5
$ws.Range("K99").Value = 5977.5
$ws.Range("M99").Value = -3731.5
$ws.Range("H132").Value = 41793370
$ws.Range("J132").Value = 41793370
$ws.Range("L132").Value = 376140330
$ws.Range("N132").Value = -376145390
$ws.Range("H134").Value = 3835.2307
$ws.Range("I134").Value = 2465.4285
$ws.Range("J134").Value = 5433.3335
$ws.Range("K134").Value = 7396.2855
$ws.Range("L134").Value = 16300.0005
$ws.Range("M134").Value = -2326.2855
$ws.Range("N134").Value = -26440.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 50000
$ws.Range("J32").Value = 50000
$ws.Range("L32").Value = 50000
$ws.Range("N32").Value = -50592
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H51").Value = 150000
$ws.Range("J51").Value = 150000
$ws.Range("L51").Value = 150000
$ws.Range("N51").Value = -151018
$ws.Range("H70").Value = 7796
$ws.Range("I70").Value = 8674.2
$ws.Range("J70").Value = 5600.5
$ws.Range("K70").Value = 8674.2
$ws.Range("L70").Value = 5600.5
$ws.Range("M70").Value = -8404.2
$ws.Range("N70").Value = -6140.5
$ws.Range("H73").Value = 7796
$ws.Range("I73").Value = 8674.2
$ws.Range("J73").Value = 5600.5
$ws.Range("K73").Value = 8674.2
$ws.Range("L73").Value = 5600.5
$ws.Range("M73").Value = -7738.200000000001
$ws.Range("N73").Value = -7472.5
$ws.Range("H80").Value = 4771.6294
$ws.Range("I80").Value = 6161.385
$ws.Range("J80").Value = 3481.1428
$ws.Range("K80").Value = 6161.385
$ws.Range("L80").Value = 3481.1428
$ws.Range("M80").Value = -5163.385
$ws.Range("N80").Value = -5477.1428
$ws.Range("H83").Value = 4771.6294
$ws.Range("I83").Value = 6161.385
$ws.Range("J83").Value = 3481.1428
$ws.Range("K83").Value = 30806.925
$ws.Range("L83").Value = 17405.714
$ws.Range("M83").Value = -25814.925
$ws.Range("N83").Value = -27389.714
$ws.Range("H102").Value = 12890.692
$ws.Range("I102").Value = 12890.692
$ws.Range("K102").Value = 12890.692
$ws.Range("M102").Value = -11268.692
$ws.Range("H128").Value = 78767.62
$ws.Range("J128").Value = 78767.62
$ws.Range("L128").Value = 78767.62
$ws.Range("N128").Value = -88727.62
$ws.Range("H132").Value = 7490.88
$ws.Range("I132").Value = 4908.7896
$ws.Range("J132").Value = 15667.5
$ws.Range("K132").Value = 14726.3688
$ws.Range("L132").Value = 47002.5
$ws.Range("M132").Value = -12196.3688
$ws.Range("N132").Value = -52062.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 98387.375
$ws.Range("I40").Value = 142521.8
$ws.Range("K40").Value = 142521.8
$ws.Range("M40").Value = -142385.8
$ws.Range("H100").Value = 4442.4614
$ws.Range("I100").Value = 4521
$ws.Range("K100").Value = 4521
$ws.Range("M100").Value = -3980
$ws.Range("H122").Value = 4900
$ws.Range("I122").Value = 3420
$ws.Range("J122").Value = 5825
$ws.Range("K122").Value = 10260
$ws.Range("L122").Value = 17475
$ws.Range("M122").Value = -7810
$ws.Range("N122").Value = -22375
$ws.Range("H132").Value = 53499.75
$ws.Range("J132").Value = 12000
$ws.Range("L132").Value = 36000
$ws.Range("N132").Value = -41060
$ws.Range("H136").Value = 4379.5835
$ws.Range("I136").Value = 2264.375
$ws.Range("J136").Value = 5437.1875
$ws.Range("K136").Value = 6793.125
$ws.Range("L136").Value = 16311.5625
$ws.Range("M136").Value = -4243.125
$ws.Range("N136").Value = -21411.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8355.111
$ws.Range("I122").Value = 6139.7
$ws.Range("J122").Value = 11124.375
$ws.Range("K122").Value = 18419.1
$ws.Range("L122").Value = 33373.125
$ws.Range("M122").Value = -15969.1
$ws.Range("N122").Value = -38273.125
$ws.Range("H123").Value = 49000
$ws.Range("J123").Value = 49000
$ws.Range("L123").Value = 49000
$ws.Range("N123").Value = -58800
$ws.Range("H126").Value = 23586.191
$ws.Range("I126").Value = 46659.332
$ws.Range("K126").Value = 139977.996
$ws.Range("M126").Value = -137507.996
$ws.Range("H132").Value = 8561
$ws.Range("I132").Value = 8830.23
$ws.Range("K132").Value = 26490.69
$ws.Range("M132").Value = -23960.69
$ws.Range("H136").Value = 4747.4116
$ws.Range("I136").Value = 3516.7693
$ws.Range("K136").Value = 10550.3079
$ws.Range("M136").Value = -8000.3079
